# CMS-research.xlsx : add product URLs to the "Name" column (column A)
# for concrete5, Liferay, and WordPress, matching the already-updated
# rows (Drupal, DotNetNuke) that already had "<name> - <url>" text.
#
# Shared-string table order in the target workbook shows the three new
# strings were appended in this order: WordPress, concrete5, Liferay —
# so we write the cells in that same order to keep the sharedStrings.xml
# layout identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRM list")

$ws.Range("A6").Value = "WordPress - https://en.wordpress.com/features/"
$ws.Range("A2").Value = "concrete5 - https://www.concrete5.org/"
$ws.Range("A5").Value = "Liferay Portal (ENTERPRISE SOLUTION ONLY) - https://www.liferay.com/product/features"

# Mirror the author's final view/selection state (scrolled down a bit,
# landed on A6 after editing that row last).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("A6").Select()
